$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: each existing dated pair of rows (138..193) shifts down to the
# position occupied by the prior pair, newest week's data lands in 138/139, and
# the oldest pair that falls off the bottom is appended as new rows 194/195.
$updates = @(
    @{Row=138; D=44455; J=800; M=650; P=130},
    @{Row=139; D=44455; J=400; M=500; P=100},
    @{Row=140; D=44427; J=600; M=650; P=130},
    @{Row=141; D=44427; J=300; M=500; P=100},
    @{Row=142; D=44343; J=1000; M=650; P=130},
    @{Row=143; D=44343; J=500; M=500; P=100},
    @{Row=144; D=44280; J=800; M=650; P=130},
    @{Row=145; D=44280; J=400; M=500; P=100},
    @{Row=146; D=44390; J=600; M=650; P=130},
    @{Row=147; D=44390; J=300; M=500; P=100},
    @{Row=148; D=44386; J=800; M=650; P=130},
    @{Row=149; D=44386; J=400; M=500; P=100},
    @{Row=150; D=44308; J=600; M=650; P=130},
    @{Row=151; D=44308; J=300; M=500; P=100},
    @{Row=152; D=44264; J=600; M=650; P=130},
    @{Row=153; D=44264; J=300; M=500; P=100},
    @{Row=154; D=44196; J=800; M=650; P=130},
    @{Row=155; D=44196; J=400; M=500; P=100},
    @{Row=156; D=44243; J=800; M=650; P=130},
    @{Row=157; D=44243; J=400; M=500; P=100},
    @{Row=158; D=44252; J=800; M=650; P=130},
    @{Row=159; D=44252; J=400; M=500; P=100},
    @{Row=160; D=44166; J=600; M=650; P=130},
    @{Row=161; D=44166; J=300; M=500; P=100},
    @{Row=162; D=44168; J=600; M=650; P=130},
    @{Row=163; D=44168; J=300; M=500; P=100},
    @{Row=164; D=44316; J=1000; M=650; P=130},
    @{Row=165; D=44316; J=500; M=500; P=100},
    @{Row=166; D=44397; J=600; M=650; P=130},
    @{Row=167; D=44397; J=300; M=500; P=100},
    @{Row=168; D=44273; J=600; M=650; P=130},
    @{Row=169; D=44273; J=300; M=500; P=100},
    @{Row=170; D=44372; J=600; M=650; P=130},
    @{Row=171; D=44372; J=300; M=500; P=100},
    @{Row=172; D=44365; J=600; M=650; P=130},
    @{Row=173; D=44365; J=300; M=500; P=100},
    @{Row=174; D=44306; J=600; M=650; P=130},
    @{Row=175; D=44306; J=300; M=500; P=100},
    @{Row=176; D=44215; J=800; M=650; P=130},
    @{Row=177; D=44215; J=400; M=500; P=100},
    @{Row=178; D=44357; J=800; M=650; P=130},
    @{Row=179; D=44357; J=400; M=500; P=100},
    @{Row=180; D=44162; J=800; M=650; P=130},
    @{Row=181; D=44162; J=400; M=500; P=100},
    @{Row=182; D=44239; J=700; M=643; P=129},
    @{Row=183; D=44239; J=300; M=500; P=100},
    @{Row=184; D=44376; J=600; M=650; P=130},
    @{Row=185; D=44376; J=300; M=500; P=100},
    @{Row=186; D=44292; J=600; M=650; P=130},
    @{Row=187; D=44292; J=300; M=500; P=100},
    @{Row=188; D=44358; J=600; M=650; P=130},
    @{Row=189; D=44358; J=300; M=500; P=100},
    @{Row=190; D=44211; J=600; M=650; P=130},
    @{Row=191; D=44211; J=300; M=500; P=100},
    @{Row=192; D=44425; J=600; M=650; P=130},
    @{Row=193; D=44425; J=300; M=500; P=100}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 10).Value = $u.J
    $ws.Cells.Item($r, 13).Value = $u.M
    $ws.Cells.Item($r, 16).Value = $u.P
}

# Append the two rows that fell off the bottom of the shifted window (old
# rows 192/193 data) as new rows 194/195, copying the static columns from the
# row directly above and the date format from an existing date cell.
$ws.Cells.Item(194, 1).Value = 11
$ws.Cells.Item(194, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(194, 3).Value = "Bíobío"
$ws.Cells.Item(194, 4).Value = 44323
$ws.Cells.Item(194, 4).NumberFormat = $ws.Range("D193").NumberFormat
$ws.Cells.Item(194, 5).Value = 8
$ws.Cells.Item(194, 6).Value = 100114014
$ws.Cells.Item(194, 7).Value = "Betarraga"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 600
$ws.Cells.Item(194, 11).Value = 600
$ws.Cells.Item(194, 12).Value = 700
$ws.Cells.Item(194, 13).Value = 650
$ws.Cells.Item(194, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(194, 15).Value = "Región Metropolitana"
$ws.Cells.Item(194, 16).Value = 130
$ws.Cells.Item(194, 17).Value = 5
$ws.Cells.Item(194, 18).Value = "Hortaliza"

$ws.Cells.Item(195, 1).Value = 11
$ws.Cells.Item(195, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(195, 3).Value = "Bíobío"
$ws.Cells.Item(195, 4).Value = 44323
$ws.Cells.Item(195, 4).NumberFormat = $ws.Range("D193").NumberFormat
$ws.Cells.Item(195, 5).Value = 8
$ws.Cells.Item(195, 6).Value = 100114014
$ws.Cells.Item(195, 7).Value = "Betarraga"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Segunda"
$ws.Cells.Item(195, 10).Value = 300
$ws.Cells.Item(195, 11).Value = 500
$ws.Cells.Item(195, 12).Value = 500
$ws.Cells.Item(195, 13).Value = 500
$ws.Cells.Item(195, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(195, 15).Value = "Región Metropolitana"
$ws.Cells.Item(195, 16).Value = 100
$ws.Cells.Item(195, 17).Value = 5
$ws.Cells.Item(195, 18).Value = "Hortaliza"

Write-Output "Updated rows 138-193 and appended rows 194-195"
